$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date serial values in column E (rows 3-5)
$ws.Range("E3").Value = (Get-Date -Year 1930 -Month 4 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E4").Value = (Get-Date -Year 1950 -Month 5 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E5").Value = (Get-Date -Year 2012 -Month 12 -Day 20 -Hour 0 -Minute 0 -Second 0)

# Autofit column E to match the new content (bestFit width)
$ws.Columns.Item(5).AutoFit() | Out-Null

# Move the active selection from C5 to E5
$ws.Range("E5").Select() | Out-Null
